$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fall 2022 block (rows 4-8) ---
# Row 4: A4 text change (FINC 1100 -> POLS 1101); E4/F4 (CPSC 4899 / 3) removed
$ws.Range("A4").Value = "POLS 1101"
$ws.Range("E4:F4").ClearContents()

# Row 5: A5 text change (POLS 1101 -> GEOL 1121K), B5 3 -> 4; C5 text change (CPSC 3415 -> CPSC 4135), D5 1 -> 3
$ws.Range("A5").Value = "GEOL 1121K"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "CPSC 4135"
$ws.Range("D5").Value = 3

# Row 6: A6 text change (GEOL 1110 -> DSCI 3111); C6 text change (CPSC 4135 -> CPSC 4148)
$ws.Range("A6").Value = "DSCI 3111"
$ws.Range("C6").Value = "CPSC 4148"

# Row 7: A7 text change (DSCI 3111 -> CPSC 3121); C7 text change (CYBR 4145 -> CPSC 4155)
$ws.Range("A7").Value = "CPSC 3121"
$ws.Range("C7").Value = "CPSC 4155"

# Row 8: A8 text change (CPSC 3121 -> CPSC 4000), B8 3 -> 0; C8 text change (CPSC 4148 -> CPSC 4157)
$ws.Range("A8").Value = "CPSC 4000"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = "CPSC 4157"

# Row 9: removed entirely (CPSC 4000 / 0)
$ws.Range("A9:B9").ClearContents()

# --- Fall 2023 block (rows 13-15) ---
# Row 13: A13 text change (CPSC 4155 -> CPSC 4175)
$ws.Range("A13").Value = "CPSC 4175"

# Row 14: A14 text change (CPSC 4157 -> CPSC 4205)
$ws.Range("A14").Value = "CPSC 4205"

# Row 15: A15 text change (CPSC 4175 -> CPSC 4555)
$ws.Range("A15").Value = "CPSC 4555"
